$d = $word.ActiveDocument

# wdFindWrap: 1 = wdFindContinue
$wdFindContinue = 1
# wdReplace: 2 = wdReplaceAll
$wdReplaceAll = 2

function Replace-Text($findText, $replaceText) {
    # MatchWholeWord is left off (not meaningful for Japanese text and can
    # cause the search to spuriously fail to match CJK strings).
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll) | Out-Null
}

# 1) Map tutorial text: rewrite the description of the town map.
Replace-Text "? 町の全体地図です。いつでも携帯メニューからアクセスできます。行きたい場所をセレクトしてみましょう。" "これが町の全体図だ。いつでも携帯メニューからアクセスできる。行きたい場所をセレクトして "

# 2) Phone tutorial "yes" branch: add a closing sentence.
Replace-Text "はい：よし、手間が省けたな。" "はい：よし、手間が省けたな。するだけだ。"

# 3) Phone tutorial body: "この町では皆" -> "ここでは皆"
Replace-Text "何でもできる。ゲームしたり、友達にメッセージしたり、世界情勢を知ることだってできる。この町では皆" "何でもできる。ゲームしたり、友達にメッセージしたり、世界情勢を知ることだってできる。ここでは皆"

# 4) Energy/stamina explanation: clarify the in-game clock wording.
Replace-Text "いくら体力が残ってるかはゲーム内時計をみてもらえれば分かる。時間が過ぎれば体力も減っていく。また誰かと会ったりすると一定量の体力が削られる。" "いくら体力が残ってるかはゲーム内の時計を見れば分かる。時間が過ぎれば体力も減っていく。また誰かと会ったりすると一定量の体力が削られる。"

# 5) Closing paragraph: rework the ending sentence.
Replace-Text "家に帰り寝ると次の日が始まる。遅くまで起きている事もできるが、寝過ごしたり学校を休んだりしてしまうこともあるので、するなら自己責任でな。どうでもいいけど。まぁやってみれば分かる。" "家に帰り寝ると次の日が始まる。遅くまで起きている事もできるが、寝過ごしたり学校を休んだりしてしまうこともあるので、するなら自己責任でな。まぁやりたければやってみればいいさ。"
